$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colBValues = @(25.09075226886648, 24.73262991661847, 24.51880584599041, 24.43330341903479, 24.41920748390457, 24.51764598043233, 24.96607864961744, 25.88869110749755, 26.58623005109778, 26.90641408821469, 27.02795507243159, 27.00176781923557, 26.91640811064536, 26.86415766999609, 26.56535284170882, 26.38269582127745, 26.27791344506329, 26.24248717550743, 26.40211215903684, 26.94147327209923, 27.29563901242966, 27.1065004502764, 26.39333330863638, 25.6351560443835)
for ($i = 0; $i -lt $colBValues.Length; $i++) {
    $ws.Cells.Item(2 + $i, 2).Value = $colBValues[$i]
}

$colCValues = @(24.15618532984284, 23.77304961131124, 23.5430714166033, 23.45079653232052, 23.43556491179757, 23.5418209681981, 24.02305841224911, 25.0032596905763, 25.73849406880737, 26.07473921807656, 26.20220031093351, 26.17474533373839, 26.08522356393809, 26.03040251113866, 25.71654485411842, 25.52437071551807, 25.4140123194585, 25.37668077241218, 25.54481077026481, 26.11151565632311, 26.48260673349757, 26.28452286654845, 25.5355694224748, 24.73492313446009)
for ($i = 0; $i -lt $colCValues.Length; $i++) {
    $ws.Cells.Item(2 + $i, 3).Value = $colCValues[$i]
}

$colDValues = @(13.66559040069222, 13.66778079066164, 13.67195642976651, 13.67436846025842, 13.67481183765661, 13.67198608480203, 13.66575737696229, 13.67605515973585, 13.69740068237532, 13.71010908626859, 13.71535239923475, 13.71420399993094, 13.71053182303065, 13.70833861439628, 13.69663047849593, 13.69021577730628, 13.68680845183042, 13.68570326653841, 13.69086942291092, 13.71159873887621, 13.72765791885568, 13.71885720032209, 13.69057303580182, 13.67085224306525)
for ($i = 0; $i -lt $colDValues.Length; $i++) {
    $ws.Cells.Item(2 + $i, 4).Value = $colDValues[$i]
}

$colEValues = @(13.71997906063499, 13.74538545379302, 13.76313299624102, 13.77090456455566, 13.77222757562498, 13.7632356237983, 13.72829290674806, 13.67684918752028, 13.64951447843202, 13.63936065491454, 13.63584437935504, 13.63658703884651, 13.6390647757698, 13.64062529717003, 13.65022401850487, 13.65669727110739, 13.66063515805525, 13.66200529989806, 13.65598596184385, 13.63832807607792, 13.6287042012195, 13.63366503437585, 13.65630687100763, 13.68893242629435)
for ($i = 0; $i -lt $colEValues.Length; $i++) {
    $ws.Cells.Item(2 + $i, 5).Value = $colEValues[$i]
}

$colGValues = @(3.783761299183795, 3.78924852137205, 3.792785553137159, 3.794269330155128, 3.79451827736247, 3.792805391899981, 3.78561857544896, 3.772847996915092, 3.764259057219844, 3.760521297853384, 3.759130048902604, 3.759428607933388, 3.760406355916358, 3.761008395139223, 3.764506720191204, 3.766696073291643, 3.767971287232394, 3.768405799262622, 3.766461362977553, 3.760118513542147, 3.756113821509069, 3.758238390395558, 3.766567424021326, 3.776162487358329)
for ($i = 0; $i -lt $colGValues.Length; $i++) {
    $ws.Cells.Item(2 + $i, 7).Value = $colGValues[$i]
}

$colIValues = @(37.30960619376872, 37.19449987793202, 37.13171405743086, 37.10811675089968, 37.10431862018095, 37.13138776048642, 37.26828325375367, 37.59905768930605, 37.87957558427436, 38.01518554264639, 38.0676721684714, 38.05631807756053, 38.01948102116788, 37.99706443089916, 37.87087279079245, 37.79549643752049, 37.75289614816707, 37.73860242623716, 37.80344243590848, 38.03027031928796, 38.18511681092302, 38.1018742072548, 37.79984776075877, 37.5029453295924)
for ($i = 0; $i -lt $colIValues.Length; $i++) {
    $ws.Cells.Item(2 + $i, 9).Value = $colIValues[$i]
}

$colJValues = @(8.27292902613156, 8.28870908821834, 8.298911861138279, 8.303199302029826, 8.30391907944473, 8.298969157078943, 8.278263689054457, 8.241712074049849, 8.217293232944083, 8.206705963337075, 8.202771181852029, 8.203615307451525, 8.206380758685963, 8.208084349194825, 8.21799557420381, 8.224208849615144, 8.227831630378198, 8.229066685472054, 8.223542361212246, 8.205566463976114, 8.194251519174767, 8.200251042004721, 8.223843522703467, 8.251170047411094)
for ($i = 0; $i -lt $colJValues.Length; $i++) {
    $ws.Cells.Item(2 + $i, 10).Value = $colJValues[$i]
}

$colMValues = @(25.83956220506195, 25.73737220286801, 25.68141914200134, 25.66033682917141, 25.6569402312358, 25.68112784529954, 25.80292474664657, 26.09507140126161, 26.34127923694151, 26.45989978400499, 26.5057468180483, 26.49583196360096, 26.46365319118933, 26.44406288046636, 26.33365830515015, 26.26760764678901, 26.23024164562415, 26.21769813137047, 26.27457437660315, 26.47307988943945, 26.60821160673213, 26.53560376756219, 26.27142282434418, 26.01041884593843)
for ($i = 0; $i -lt $colMValues.Length; $i++) {
    $ws.Cells.Item(2 + $i, 13).Value = $colMValues[$i]
}

$colNValues = @(22.80963733486907, 22.83731474705726, 22.85620908691634, 22.86438491372342, 22.86577121944257, 22.85631742284983, 22.81878512513221, 22.76033595268069, 22.72673295718, 22.71349656441296, 22.70878090942264, 22.70978328431359, 22.71310264766772, 22.7151745473787, 22.72763939743285, 22.73581254011855, 22.74070633449901, 22.74239634757388, 22.73492252523445, 22.71211960295649, 22.69894675610074, 22.70581840045426, 22.73532429389666, 22.77451535718361)
for ($i = 0; $i -lt $colNValues.Length; $i++) {
    $ws.Cells.Item(2 + $i, 14).Value = $colNValues[$i]
}

